$d = $word.ActiveDocument

# Move to the very end of the document body (end of the last paragraph)
$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)  # wdCollapseEnd

# Paragraph: "Install .Net 6 runtime" (top-level bullet, ilvl 0)
$endRange.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "Install .Net 6 runtime"
$p1.Range.ListFormat.ListLevelNumber = 1

# Paragraph: "sudo apt-get install -y aspnetcore-runtime-6.0 " (sub bullet, ilvl 1)
$endRange2 = $p1.Range
$endRange2.Collapse(0)
$endRange2.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "sudo apt-get install -y aspnetcore-runtime-6.0 "
$p2.Range.ListFormat.ListLevelNumber = 2

# Paragraph: empty bullet (top-level, ilvl 0)
$endRange3 = $p2.Range
$endRange3.Collapse(0)
$endRange3.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.ListFormat.ListLevelNumber = 1
